# Lotka_Volterra_realistic.xlsx - "Add last tweaks to the code before big launch"
#
# The only content-level change in this commit is a rewrite of the reaction-term
# formula text stored in the "Reaction terms" row (cell E2 on Sheet1): the
# en-dash minus sign is replaced with a plain hyphen-minus and the first
# multiplicative term gets an extra pair of parentheses around it:
#   N*( r*(1-(N/K)) - ((k*P)/(N+D)))   ->   N*((r*(1-(N/K)))-((k*P)/(N+D)))
#
# (Everything else in the upstream diff - the workbookPr attribute spelling,
# the styles.xml font table layout, and the sheet's default column width -
# is a cosmetic artifact of the authoring tool re-serialising the file and
# carries no semantic content, so it isn't reproduced here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "N*((r*(1-(N/K)))-((k*P)/(N+D)))"

# Workbook date system is unaffected by this edit (still 1900 date system);
# set explicitly for clarity / parity with the workbookPr flag being touched.
$wb.Date1904 = $false
